# Auto-generated edit script: updates cryptocurrency Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.276.94"
$ws.Range("E2").Value = "  -0.35%  "
$ws.Range("D3").Value = "3.137.40"
$ws.Range("E3").Value = "  -1.41%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "163.83"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.32%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.574"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.55%  "
$ws.Range("D9").Value = "3.150.25"
$ws.Range("E9").Value = "  -1.29%  "
$ws.Range("E10").Value = "  -3.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.60"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.97%  "
$ws.Range("E12").Value = "  +0.07%  "
$ws.Range("D13").Value = "3.682.77"
$ws.Range("E13").Value = "  -1.52%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.127"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.34%  "
$ws.Range("D15").Value = "64.313.88"
$ws.Range("E15").Value = "  -0.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.11"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.28%  "
$ws.Range("D17").Value = "3.138.97"
$ws.Range("E17").Value = "  -1.56%  "
$ws.Range("E18").Value = "  -3.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "401.92"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.25%  "
$ws.Range("E20").Value = "  -1.86%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.08"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.81%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.69"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.67%  "
$ws.Range("E25").Value = "  -0.81%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.197"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.81%  "
$ws.Range("E27").Value = "  -4.82%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.78"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.996"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.30%  "
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("E31").Value = "  -1.62%  "
$ws.Range("E32").Value = "  -2.60%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "161.57"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.41%  "
$ws.Range("E34").Value = "  -0.94%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.85"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.91%  "
$ws.Range("E36").Value = "  -2.56%  "
$ws.Range("E37").Value = "  -1.90%  "
$ws.Range("E38").Value = "  -1.99%  "
$ws.Range("D39").Value = "2.640.07"
$ws.Range("E39").Value = "  -3.52%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "23.75"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.06"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.38%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "38.37"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.690"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.68%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0615"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.42"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.52%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0254"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.89%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "21.15"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.83%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "286.36"
$ws.Range("D48").Style = "Normal"
$ws.Range("E49").Value = "  -0.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0975"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "10.46"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.06%  "
